$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G8").Value = 0.3921
$ws.Range("F9").Value = 0.8801
$ws.Range("G9").Value = 0.6433
$ws.Range("G10").Value = 0.346
